# Weekly update: insert this week's new Coliflor price rows (Vega Central
# Mapocho de Santiago) right after the existing row 314, shifting the rest
# of the table down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new blank rows starting at row 315 (each Insert() pushes the
# current row 315 and everything below it down by one).
$ws.Rows.Item(315).Insert()
$ws.Rows.Item(315).Insert()
$ws.Rows.Item(315).Insert()
$ws.Rows.Item(315).Insert()

# Row 315: Primera, Región Metropolitana
$ws.Cells.Item(315, 1).Value = 9
$ws.Cells.Item(315, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(315, 3).Value = "Metropolitana"
$ws.Cells.Item(315, 4).Value = "2021-09-09"
$ws.Cells.Item(315, 5).Value = 13
$ws.Cells.Item(315, 6).Value = 100112008
$ws.Cells.Item(315, 7).Value = "Coliflor"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 3400
$ws.Cells.Item(315, 11).Value = 600
$ws.Cells.Item(315, 12).Value = 650
$ws.Cells.Item(315, 13).Value = 625
$ws.Cells.Item(315, 14).Value = "`$/unidad"
$ws.Cells.Item(315, 15).Value = "Región Metropolitana"
$ws.Cells.Item(315, 16).Value = 625
$ws.Cells.Item(315, 17).Value = 1
$ws.Cells.Item(315, 18).Value = "Hortaliza"

# Row 316: Primera, Región de O'Higgins
$ws.Cells.Item(316, 1).Value = 9
$ws.Cells.Item(316, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(316, 3).Value = "Metropolitana"
$ws.Cells.Item(316, 4).Value = "2021-09-09"
$ws.Cells.Item(316, 5).Value = 13
$ws.Cells.Item(316, 6).Value = 100112008
$ws.Cells.Item(316, 7).Value = "Coliflor"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 4300
$ws.Cells.Item(316, 11).Value = 600
$ws.Cells.Item(316, 12).Value = 650
$ws.Cells.Item(316, 13).Value = 625
$ws.Cells.Item(316, 14).Value = "`$/unidad"
$ws.Cells.Item(316, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(316, 16).Value = 625
$ws.Cells.Item(316, 17).Value = 1
$ws.Cells.Item(316, 18).Value = "Hortaliza"

# Row 317: Segunda, Región Metropolitana
$ws.Cells.Item(317, 1).Value = 9
$ws.Cells.Item(317, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(317, 3).Value = "Metropolitana"
$ws.Cells.Item(317, 4).Value = "2021-09-09"
$ws.Cells.Item(317, 5).Value = 13
$ws.Cells.Item(317, 6).Value = 100112008
$ws.Cells.Item(317, 7).Value = "Coliflor"
$ws.Cells.Item(317, 8).Value = "Sin especificar"
$ws.Cells.Item(317, 9).Value = "Segunda"
$ws.Cells.Item(317, 10).Value = 970
$ws.Cells.Item(317, 11).Value = 450
$ws.Cells.Item(317, 12).Value = 500
$ws.Cells.Item(317, 13).Value = 475
$ws.Cells.Item(317, 14).Value = "`$/unidad"
$ws.Cells.Item(317, 15).Value = "Región Metropolitana"
$ws.Cells.Item(317, 16).Value = 475
$ws.Cells.Item(317, 17).Value = 1
$ws.Cells.Item(317, 18).Value = "Hortaliza"

# Row 318: Segunda, Región de O'Higgins
$ws.Cells.Item(318, 1).Value = 9
$ws.Cells.Item(318, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(318, 3).Value = "Metropolitana"
$ws.Cells.Item(318, 4).Value = "2021-09-09"
$ws.Cells.Item(318, 5).Value = 13
$ws.Cells.Item(318, 6).Value = 100112008
$ws.Cells.Item(318, 7).Value = "Coliflor"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Segunda"
$ws.Cells.Item(318, 10).Value = 1600
$ws.Cells.Item(318, 11).Value = 450
$ws.Cells.Item(318, 12).Value = 500
$ws.Cells.Item(318, 13).Value = 475
$ws.Cells.Item(318, 14).Value = "`$/unidad"
$ws.Cells.Item(318, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(318, 16).Value = 475
$ws.Cells.Item(318, 17).Value = 1
$ws.Cells.Item(318, 18).Value = "Hortaliza"
